# prepare for human readabily
# Strip the stray trailing "16" that was accidentally appended to every
# scripture reference in column A (e.g. "Nahum 1:216" -> "Nahum 1:2").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $value = $cell.Value()
    if ($value -ne $null -and $value.ToString().EndsWith("16")) {
        $cell.Value = $value.ToString().Substring(0, $value.ToString().Length - 2)
    }
}
